$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text {
    param($addr, $val)
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

# Row 2-8: simple price updates
Set-Text "D2" "244.63"
Set-Text "D3" "21.82"
Set-Text "D4" "5.402"
Set-Text "D5" "0.06044"
Set-Text "D6" "3.398"
Set-Text "D7" "0.8149"
Set-Text "D8" "0.9240"

# Rows 9-17: coin list rotates by one position with updated prices
Set-Text "B9" "WazirX"
Set-Text "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-Text "D9" "0.1441"
Set-Text "E9" "8WazirXWRX"

Set-Text "B10" "MandalaExchangeToken"
Set-Text "C10" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-Text "D10" "0.07475"
Set-Text "E10" "9MandalaExchangeTokenMDX"

Set-Text "B11" "LiechtensteinCryptoassetsExchange"
Set-Text "C11" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-Text "D11" "0.03382"
Set-Text "E11" "10LiechtensteinCryptoassetsExchangeLCX"

Set-Text "B12" "BitrueCoin"
Set-Text "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-Text "D12" "0.03046"
Set-Text "E12" "11BitrueCoinBTR"

Set-Text "B13" "BitMartToken"
Set-Text "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-Text "D13" "0.09404"
Set-Text "E13" "12BitMartTokenBMX"

Set-Text "B14" "MCDex"
Set-Text "C14" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-Text "D14" "4.005"
Set-Text "E14" "13MCDexMCB"

Set-Text "B15" "BitForexToken"
Set-Text "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-Text "D15" "0.001589"
Set-Text "E15" "14BitForexTokenBF"

Set-Text "B16" "CoinExToken"
Set-Text "C16" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-Text "D16" "0.04795"
Set-Text "E16" "15CoinExTokenCET"

Set-Text "B17" "One"
Set-Text "C17" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-Text "D17" "0.0005944"
Set-Text "E17" "16OneONE"

# Rows 18-27: simple price updates
Set-Text "D18" "0.005485"
Set-Text "D19" "0.004161"
Set-Text "D20" "0.0009857"
Set-Text "D21" "0.00008798"
Set-Text "D23" "6.428"
Set-Text "D26" "0.1323"
Set-Text "D27" "0.0002902"

# Rows 41-43: coin list rotates by one position with updated prices
Set-Text "B41" "BKEXToken"
Set-Text "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-Text "D41" "0.1075"
Set-Text "E41" "40BKEXTokenBKK"

Set-Text "B42" "CEJI"
Set-Text "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-Text "D42" "0.002719"
Set-Text "E42" "41CEJICEJI"

Set-Text "B43" "KickToken"
Set-Text "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-Text "D43" "0.003054"
Set-Text "E43" "42KickTokenKICK"

# Rows 44-50: simple price updates
Set-Text "D44" "0.006376"
Set-Text "D45" "0.00005244"
Set-Text "D46" "0.00000000751"
Set-Text "D47" "1.101"
Set-Text "E47" "46CoinbaseStockTokenCOINBestin24h"
Set-Text "D48" "0.002528"
Set-Text "D49" "0.00002101"
Set-Text "D50" "0.01011"
